$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are rotated among rows: D, K, L, M, N, O, P, R, S
$cols = @("D", "K", "L", "M", "N", "O", "P", "R", "S")

# Snapshot the current ("before") values for the rows involved in the rotation
$rows = @(2, 3, 5, 6, 7, 8)
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowData
}

# New row <- old row mapping (cyclic rotation of data between rows)
$mapping = @{
    2 = 7
    3 = 6
    5 = 2
    6 = 8
    7 = 3
    8 = 5
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $snapshot[$oldRow][$c]
    }
}
